{"js": "// Update the two-digit \u00f7 one-digit division drill table: replace the\n// text of 25 table cells (5 data rows \u00d7 5 columns) with newly generated\n// problems. Old values are unique within the document and are applied\n// strictly in document order so that a freshly written \"new\" value that\n// happens to equal another cell's old value is never matched prematurely.\nconst replacements = [\n  [\"33\u00f73=\", \"90\u00f74=\"],\n  [\"88\u00f78=\", \"45\u00f74=\"],\n  [\"96\u00f78=\", \"19\u00f73=\"],\n  [\"14\u00f77=\", \"18\u00f74=\"],\n  [\"18\u00f72=\", \"76\u00f78=\"],\n  [\"80\u00f78=\", \"82\u00f78=\"],\n  [\"33\u00f72=\", \"31\u00f79=\"],\n  [\"55\u00f72=\", \"96\u00f77=\"],\n  [\"29\u00f72=\", \"16\u00f75=\"],\n  [\"79\u00f79=\", \"29\u00f72=\"],\n  [\"33\u00f75=\", \"21\u00f72=\"],\n  [\"81\u00f72=\", \"70\u00f74=\"],\n  [\"20\u00f78=\", \"81\u00f75=\"],\n  [\"56\u00f78=\", \"47\u00f74=\"],\n  [\"78\u00f76=\", \"98\u00f76=\"],\n  [\"40\u00f77=\", \"90\u00f76=\"],\n  [\"99\u00f74=\", \"19\u00f72=\"],\n  [\"46\u00f77=\", \"32\u00f79=\"],\n  [\"19\u00f78=\", \"63\u00f77=\"],\n  [\"24\u00f77=\", \"31\u00f74=\"],\n  [\"75\u00f74=\", \"83\u00f76=\"],\n  [\"84\u00f76=\", \"33\u00f72=\"],\n  [\"76\u00f72=\", \"69\u00f74=\"],\n  [\"86\u00f77=\", \"39\u00f74=\"],\n  [\"22\u00f78=\", \"20\u00f79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Exactly one occurrence is expected at this point in the sequence.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the two-digit \u00f7 one-digit division drill table: replace the\n# text of 25 table cells (5 data rows \u00d7 5 columns) with newly generated\n# problems. Old values are unique within the document and are applied\n# strictly in document order so that a freshly written \"new\" value that\n# happens to equal another cell's old value is never matched prematurely.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"33\u00f73=\", \"90\u00f74=\"),\n  @(\"88\u00f78=\", \"45\u00f74=\"),\n  @(\"96\u00f78=\", \"19\u00f73=\"),\n  @(\"14\u00f77=\", \"18\u00f74=\"),\n  @(\"18\u00f72=\", \"76\u00f78=\"),\n  @(\"80\u00f78=\", \"82\u00f78=\"),\n  @(\"33\u00f72=\", \"31\u00f79=\"),\n  @(\"55\u00f72=\", \"96\u00f77=\"),\n  @(\"29\u00f72=\", \"16\u00f75=\"),\n  @(\"79\u00f79=\", \"29\u00f72=\"),\n  @(\"33\u00f75=\", \"21\u00f72=\"),\n  @(\"81\u00f72=\", \"70\u00f74=\"),\n  @(\"20\u00f78=\", \"81\u00f75=\"),\n  @(\"56\u00f78=\", \"47\u00f74=\"),\n  @(\"78\u00f76=\", \"98\u00f76=\"),\n  @(\"40\u00f77=\", \"90\u00f76=\"),\n  @(\"99\u00f74=\", \"19\u00f72=\"),\n  @(\"46\u00f77=\", \"32\u00f79=\"),\n  @(\"19\u00f78=\", \"63\u00f77=\"),\n  @(\"24\u00f77=\", \"31\u00f74=\"),\n  @(\"75\u00f74=\", \"83\u00f76=\"),\n  @(\"84\u00f76=\", \"33\u00f72=\"),\n  @(\"76\u00f72=\", \"69\u00f74=\"),\n  @(\"86\u00f77=\", \"39\u00f74=\"),\n  @(\"22\u00f78=\", \"20\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 0\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
